$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.450817465782166
$ws.Range("B1").Value = 2.068468809127808
$ws.Range("C1").Value = 3.373941421508789
$ws.Range("D1").Value = 5.545896530151367
$ws.Range("E1").Value = 2.378536462783813
